# PlayerPerformance_4826.xlsx edit
# - add "Player Info" sheet (before "ODI Batting")
# - add "ODI Batting Extra" sheet (after "ODI Bowling")
# - on "ODI Batting" and "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE,
#   replace the howstat URL values with just the bare match-code number
# - clear the stray empty INNING_NUMBER cells on rows where the player did not bat

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new sheets first (a Worksheet handle captured before an
#    insertion is anchored to an index, not a name, so re-resolve every
#    worksheet reference by name AFTER all Add() calls are done).
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$playerInfo.Name = "Player Info"

$extra = $wb.Worksheets.Add($null, $wb.Worksheets.Item("ODI Bowling"))
$extra.Name = "ODI Batting Extra"

# Re-fetch every sheet fresh, by name, now that the final tab order is set.
$playerInfo = $wb.Worksheets.Item("Player Info")
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Item("ODI Batting Extra")

# ---------------------------------------------------------------------------
# 2. Populate "Player Info"
# ---------------------------------------------------------------------------
$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 0; $c -lt $playerInfoHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $playerInfoHeaders[$c]
    $cell.Font.Bold = $true
}

$playerInfoRow = @("4826", "Fabian Anthony Allen", "Right Handed", "Left Arm Orthodox")
for ($c = 0; $c -lt $playerInfoRow.Length; $c++) {
    $cell = $playerInfo.Cells.Item(2, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $playerInfoRow[$c]
}

# ---------------------------------------------------------------------------
# 3. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (header + values),
#    and drop the leftover empty INNING_NUMBER cells
# ---------------------------------------------------------------------------
$battingWs.Range("D1").NumberFormat = "@"
$battingWs.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4219", "4220", "4221", "4230", "4291", "4293", "4296", "4338", "4344", "4348", "4359", "4362", "4414", "4417", "4449", "4450", "4451", "4533", "4535", "4536")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    $cell = $battingWs.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$i]
}

foreach ($row in @(8, 12, 16, 18)) {
    $battingWs.Cells.Item($row, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 4. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (header + values)
# ---------------------------------------------------------------------------
$bowlingWs.Range("B1").NumberFormat = "@"
$bowlingWs.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4219", "4220", "4221", "4230", "4291", "4293", "4296", "4338", "4344", "4348", "4362", "4414", "4449", "4450", "4451", "4533", "4535", "4536")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    $cell = $bowlingWs.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$i]
}

# ---------------------------------------------------------------------------
# 5. Populate "ODI Batting Extra"
# ---------------------------------------------------------------------------
$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 0; $c -lt $extraHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $extraHeaders[$c]
    $cell.Font.Bold = $true
}

# Each row: MATCH_CODE, BATTING_POSITION (numeric or $null), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# ($null entries stay fully blank - matching the sparsely-populated source sheet)
$extraData = @(
    @("4219", 8, "1", "0", "1.77%", "NO"),
    @("4220", 8, "1", "0", "6.54%", "NO"),
    @("4221", 7, "0", "0", "3.85%", "NO"),
    @("4230", 8, "0", "1", "3.03%", "NO"),
    @("4291", 7, "0", "0", $null, "NO"),
    @("4293", 7, "1", "0", "2.83%", "NO"),
    @("4296", 7, $null, $null, $null, "NO"),
    @("4338", 8, "0", "0", $null, "NO"),
    @("4344", 8, "7", "1", "16.19%", "NO"),
    @("4348", 8, "0", "0", $null, "NO"),
    @("4359", $null, $null, $null, $null, "NO"),
    @("4362", $null, $null, $null, $null, "NO"),
    @("4414", $null, $null, $null, $null, "NO"),
    @("4417", 7, "2", "3", "12.29%", "NO"),
    @("4449", 8, $null, $null, $null, "NO"),
    @("4450", $null, $null, $null, $null, "NO"),
    @("4451", 8, $null, $null, $null, "NO"),
    @("4533", $null, $null, $null, $null, "NO"),
    @("4535", 8, "0", "1", "6.74%", "NO"),
    @("4536", 7, "0", "0", $null, "NO")
)

for ($i = 0; $i -lt $extraData.Length; $i++) {
    $row = $i + 2
    $rowData = $extraData[$i]

    $codeCell = $extra.Cells.Item($row, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $rowData[0]

    # BATTING_POSITION is stored as a real number (when present)
    if ($rowData[1] -ne $null) {
        $extra.Cells.Item($row, 2).Value = $rowData[1]
    }

    for ($c = 2; $c -lt 5; $c++) {
        $v = $rowData[$c]
        if ($v -ne $null) {
            $cell = $extra.Cells.Item($row, $c + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $v
        }
    }

    $moCell = $extra.Cells.Item($row, 6)
    $moCell.NumberFormat = "@"
    $moCell.Value = $rowData[5]
}

# ---------------------------------------------------------------------------
# 6. Restore the originally-active tab ("Player Info" is the new first sheet)
# ---------------------------------------------------------------------------
$playerInfo.Activate()
